# Update gh-pages to output generated at 456a3b4
# Updates the "F" column ("想去人数" / want-to-go headcount) values across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 19
$ws.Range("F3").Value = 8003
$ws.Range("F5").Value = 942
$ws.Range("F6").Value = 296
$ws.Range("F7").Value = 806
$ws.Range("F9").Value = 95
$ws.Range("F10").Value = 69
$ws.Range("F12").Value = 870
$ws.Range("F13").Value = 3278
$ws.Range("F14").Value = 210
$ws.Range("F15").Value = 102
$ws.Range("F16").Value = 742
$ws.Range("F17").Value = 755
$ws.Range("F21").Value = 263
$ws.Range("F22").Value = 232
$ws.Range("F23").Value = 334
$ws.Range("F24").Value = 288
$ws.Range("F25").Value = 132
$ws.Range("F26").Value = 120
$ws.Range("F27").Value = 283
$ws.Range("F28").Value = 29
$ws.Range("F32").Value = 549
$ws.Range("F33").Value = 26
$ws.Range("F34").Value = 36
$ws.Range("F35").Value = 14
$ws.Range("F36").Value = 21
$ws.Range("F37").Value = 226
$ws.Range("F38").Value = 104

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 208

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 208
$ws.Range("F4").Value = 19
$ws.Range("F5").Value = 8004
$ws.Range("F7").Value = 942
$ws.Range("F8").Value = 296
$ws.Range("F9").Value = 806
$ws.Range("F11").Value = 95
$ws.Range("F12").Value = 69
$ws.Range("F14").Value = 871
$ws.Range("F16").Value = 3278
$ws.Range("F17").Value = 210
$ws.Range("F18").Value = 102
$ws.Range("F20").Value = 742
$ws.Range("F21").Value = 755
$ws.Range("F26").Value = 263
$ws.Range("F27").Value = 232
$ws.Range("F28").Value = 334
$ws.Range("F29").Value = 288
$ws.Range("F30").Value = 132
$ws.Range("F31").Value = 120
$ws.Range("F32").Value = 283
$ws.Range("F33").Value = 29
$ws.Range("F37").Value = 549
$ws.Range("F38").Value = 26
$ws.Range("F39").Value = 36
$ws.Range("F40").Value = 14
$ws.Range("F41").Value = 21
$ws.Range("F42").Value = 226
$ws.Range("F43").Value = 104

